$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates per the cryptos list refresh.
# For cells whose new text would otherwise be auto-parsed by Excel as a
# number (losing formatting like trailing zeros, e.g. "1.000" -> 1), force
# the cell to Text format first so the exact original string is preserved.

$ws.Range('D2').Value = '26.511.37'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.839.62'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '260.70'
$ws.Range('E5').Value = '  +0.08%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5257'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3191'
$ws.Range('E8').Value = '  -1.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06796'
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.79'
$ws.Range('E10').Value = '  +0.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7851'
$ws.Range('E11').Value = '  +2.66%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07764'
$ws.Range('E12').Value = '  +1.14%  '
$ws.Range('D13').Value = '1.837.20'
$ws.Range('E13').Value = '  -0.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.93'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.017'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.87'
$ws.Range('E17').Value = '  -0.57%  '
$ws.Range('E18').Value = '  +0.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007951'
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = '26.533.56'
$ws.Range('D21').Value = '2.075.18'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.630'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.974'
$ws.Range('E23').Value = '  +0.56%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.352'
$ws.Range('E24').Value = '  -1.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '141.90'
$ws.Range('E25').Value = '  -1.97%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.185'
$ws.Range('E26').Value = '  -1.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.678'
$ws.Range('E27').Value = '  +1.74%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.91'
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '111.53'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.167'
$ws.Range('E30').Value = '  -0.32%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08702'
$ws.Range('E31').Value = '  -0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.078'
$ws.Range('E32').Value = '  -1.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04870'
$ws.Range('E33').Value = '  +1.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7303'
$ws.Range('E34').Value = '  +4.07%  '
$ws.Range('E35').Value = '  +1.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.864'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.095'
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('E38').Value = '  +2.53%  '
$ws.Range('E39').Value = '  -0.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4806'
$ws.Range('E40').Value = '  -0.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8948'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '109.66'
$ws.Range('E42').Value = '  -1.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.947'
$ws.Range('E43').Value = '  -2.48%  '
$ws.Range('E44').Value = '  +0.10%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.672'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4183'
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.967'
$ws.Range('E47').Value = '  -0.44%  '
$ws.Range('B48').Value = 'Algorand'
$ws.Range('C48').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1234'
$ws.Range('E48').Value = '  +1.08%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05850'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.92'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.8934'
$ws.Range('E51').Value = '  +1.22%  '
